$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

$ws.Cells.Item($row, 1).Value = 45866.33356522125
$ws.Cells.Item($row, 2).Value = 2025
$ws.Cells.Item($row, 3).Value = 31
$ws.Cells.Item($row, 4).Value = 14.78
$ws.Cells.Item($row, 5).Value = 91.06999999999999
$ws.Cells.Item($row, 6).Value = 172.27
$ws.Cells.Item($row, 7).Value = 5.07
$ws.Cells.Item($row, 8).Value = "W"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = "08:00:20"

# Match the date/time number format used by the rest of column A
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
